$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost/Nord (Q/R) coordinate columns to whole numbers for rows 12 and 17
# (these rows keep their species data, only Q/R get rounded).
$ws.Range("Q12").Value = 657162
$ws.Range("R12").Value = 6571271

$ws.Range("Q17").Value = 657216
$ws.Range("R17").Value = 6571313

# Rows 13-16 have their species/observation data rotated among themselves
# (row13<-old row15, row14<-old row16, row15<-old row14, row16<-old row13),
# with Q/R rounded to whole numbers in the process.

# New row 13 (was row 15's data)
$ws.Range("A13").Value = 112128573
$ws.Range("B13").Value = 90668
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 788
$ws.Range("F13").Value = "Gul taggsvamp"
$ws.Range("G13").Value = "Hydnellum geogenium"
$ws.Range("H13").Value = "(Fr.) Banker"
$ws.Range("Q13").Value = 657134
$ws.Range("R13").Value = 6571219
$ws.Range("AC13").ClearContents()

# New row 14 (was row 16's data)
$ws.Range("A14").Value = 112128627
$ws.Range("B14").Value = 90087
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 3298
$ws.Range("F14").Value = "Trådticka"
$ws.Range("G14").Value = "Climacocystis borealis"
$ws.Range("H14").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q14").Value = 657182
$ws.Range("R14").Value = 6571192
$ws.Range("AC14").Value = "På nedre delen av torrgran."

# New row 15 (was row 14's data)
$ws.Range("A15").Value = 112128498
$ws.Range("B15").Value = 90018
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 1339
$ws.Range("F15").Value = "Brandticka"
$ws.Range("G15").Value = "Pycnoporellus fulgens"
$ws.Range("H15").Value = "(Fr.) Donk"
$ws.Range("Q15").Value = 657134
$ws.Range("R15").Value = 6571271
$ws.Range("AC15").Value = "På granlåga. En del årsfärska dödade granar av granbarkborre. Gott om död ved i form av torrträd och lågor av gran."

# New row 16 (was row 13's data)
$ws.Range("A16").Value = 112128530
$ws.Range("B16").Value = 89802
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 5420
$ws.Range("F16").Value = "Grovticka"
$ws.Range("G16").Value = "Phaeolus schweinitzii"
$ws.Range("H16").Value = "(Fr.) Pat."
$ws.Range("Q16").Value = 657144
$ws.Range("R16").Value = 6571278
$ws.Range("AC16").Value = "På högstubbe av tall."

# Clear the Starttid (Z) and Sluttid (AB) columns for all affected rows (12-17)
$ws.Range("Z12:Z17").ClearContents()
$ws.Range("AB12:AB17").ClearContents()
